# Add season "record" columns (Wins / Losses / Ties) to the roster sheet.
# The sheet currently spans A1:AB41; we extend it with three new columns
# (AC, AD, AE) carrying the team's win/loss/tie record for every player row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the header formatting (bold font, thin border, centered alignment)
# from an existing header cell so the new headers match the look of the
# rest of row 1 instead of minting a brand new style.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("AC1:AE1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Fill in the season record for every player row (2 through 41) with the
# team's overall win/loss/tie totals for the season.
for ($r = 2; $r -le 41; $r++) {
    $ws.Cells.Item($r, 29).Value = 77
    $ws.Cells.Item($r, 30).Value = 85
    $ws.Cells.Item($r, 31).Value = 0
}
